$d = $word.ActiveDocument

# 1. Title (appears twice: main heading + bold recap near the end)
$d.Content.Find.Execute(
    "Play Magic Stars 6 free slot - RTP at 96.5% - Review", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Play Magic Stars 6 Free | Exciting Features and Impressive Graphics",
    2, $false, $false, $false, $false
)

# 2. "What we like" bullet list items
$d.Content.Find.Execute(
    "Land Wilds help complete winning paylines", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Packed with special features including Land Wilds and Scatter Symbols",
    2, $false, $false, $false, $false
)

$d.Content.Find.Execute(
    "Scatter Symbols trigger up to 30 free spins with x3 multiplier", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Customizable gameplay through Wazdan's unique features",
    2, $false, $false, $false, $false
)

$d.Content.Find.Execute(
    "Customizable gameplay with unique features", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Variety of features to accommodate different player preferences",
    2, $false, $false, $false, $false
)

$d.Content.Find.Execute(
    "Impressive graphics transport players to a magnificent universe", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Impressive graphics that transport you to a magnificent universe",
    2, $false, $false, $false, $false
)

# 3. "What we don't like" bullet list items
$d.Content.Find.Execute(
    "Limited in-game bonus features", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Limited options for players who prefer medium volatility",
    2, $false, $false, $false, $false
)

$d.Content.Find.Execute(
    "No progressive jackpot", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Not a wide variety of games with the same theme",
    2, $false, $false, $false, $false
)

# 4. Meta description (italic run near the end)
$d.Content.Find.Execute(
    "Read our review of Magic Stars 6, an online slot game by Wazdan. Play for free with RTP at 96.5%. Customize gameplay, impressive graphics, and more.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Play Magic Stars 6 for free and enjoy its special features and impressive graphics.",
    2, $false, $false, $false, $false
)
